# Update the "想去人数" (column F) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    3  = 1352
    4  = 154
    6  = 231
    7  = 100
    9  = 182
    10 = 133
    11 = 4562
    12 = 6830
    16 = 570
    17 = 54
    18 = 4132
    19 = 493
    20 = 74
    22 = 2708
    26 = 355
    27 = 362
    28 = 398
    30 = 38
    31 = 1625
    34 = 137
    36 = 545
    37 = 497
    40 = 100
    41 = 643
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
